$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-6 from 45208 to 45212
$ws.Range("C2:C6").Value = 45212
